$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add RSI values in column E (rows 2-6)
$ws.Range("E2").Value = 62.8
$ws.Range("E3").Value = 17.8
$ws.Range("E4").Value = 36.1
$ws.Range("E5").Value = 26.4
$ws.Range("E6").Value = 20.8

# Update 점수(룰) (rule score) and 최종점수 (final score) for row 2
$ws.Range("G2").Value = 60
$ws.Range("K2").Value = 69

# Update MACRO_SCORE values for rows 2-6
$ws.Range("N2").Value = 85.87127175646313
$ws.Range("N3").Value = 85.87127175646313
$ws.Range("N4").Value = 85.87127175646313
$ws.Range("N5").Value = 85.87127175646313
$ws.Range("N6").Value = 85.87127175646313
